$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data text for the destinations data years (Annual Population Survey rows).
# Rows 2-4 (Employment volumes / Employment by occupation / Employment by industry):
# the "Latest period" and "Next period" columns roll forward one quarter.
$ws.Range("C2").Value = "Jul 2024 - Jun 2025 (14/10/25)"
$ws.Range("D2").Value = "Oct 2024 - Sep 2025 (20/01/26)"

$ws.Range("C3").Value = "Jul 2024 - Jun 2025 (14/10/25)"
$ws.Range("D3").Value = "Oct 2024 - Sep 2025 (20/01/26)"

$ws.Range("C4").Value = "Jul 2024 - Jun 2025 (14/10/25)"
$ws.Range("D4").Value = "Oct 2024 - Sep 2025 (20/01/26)"

# Update the active cell selection to reflect where the author was last working.
$ws.Range("D4").Select()
